# This script transforms the LOQ4241 worksheet from its original content
# (before.xlsx) into the edited target state described by the commit diff.
# The edit removes two long descriptive paragraphs (the "Programa resumido"
# summary text and the detailed numbered "Programa" list). This shifts every
# subsequent row's label up by one position, while some of the B/C data
# values get reassigned along the way, and the final row (which held the
# "Requisitos" detail text) disappears once everything shifts up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# --- Row 10: "Objetivos:" value changes to the professor's name/id ---
$ws.Range("B10").Value = "5840917 - Fabrício Maciel Gomes"
$ws.Range("C10").Value = "5840917 - Fabrício Maciel Gomes"

# --- Row 13: becomes "Programa resumido:" / "Semestral" ---
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# --- Row 14: becomes just "Short syllabus:" (B/C cleared entirely) ---
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14:C14").Clear()

# --- Row 15: becomes "Programa:" / "01/01/2016" (kept as text, not a date) ---
$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "01/01/2016"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "01/01/2016"
# restore the normal (non-text-forced) cell formatting used elsewhere in the sheet
$ws.Range("B19").Copy()
$ws.Range("B15").PasteSpecial($xlPasteFormats)
$ws.Range("C19").Copy()
$ws.Range("C15").PasteSpecial($xlPasteFormats)

# --- Row 16: becomes just "Syllabus:" (B/C cleared entirely) ---
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16:C16").Clear()

# --- Row 17: becomes just "Avaliação:" ---
$ws.Range("A17").Value = "Avaliação:"

# --- Row 18: becomes "Método:" / professor's name/id (B/C newly populated) ---
$ws.Range("A18").Value = "Método:"
$ws.Range("B19").Copy()
$ws.Range("B18").PasteSpecial($xlPasteFormats)
$ws.Range("C19").Copy()
$ws.Range("C18").PasteSpecial($xlPasteFormats)
$ws.Range("B18").Value = "5840917 - Fabrício Maciel Gomes"
$ws.Range("C18").Value = "5840917 - Fabrício Maciel Gomes"

# --- Row 19: becomes "Critério:" / (teaching method text stays) ---
$ws.Range("A19").Value = "Critério:"
$ws.Range("B19").Value = "Aulas expositivas teóricas, aulas práticas, aulas de exercícios."
$ws.Range("C19").Value = "Aulas expositivas teóricas, aulas práticas, aulas de exercícios."

# --- Row 20: becomes "Norma de recuperação:" / (grading formula stays) ---
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = "A Nota Final do aluno será determinada segundo a seguinte equação: Nota Final = (Prova- Bimestral-1*0,4) + (Prova-Bimestral-2*0,4) + (Trabalho*0,2)"
$ws.Range("C20").Value = "A Nota Final do aluno será determinada segundo a seguinte equação: Nota Final = (Prova- Bimestral-1*0,4) + (Prova-Bimestral-2*0,4) + (Trabalho*0,2)"

# --- Row 21: becomes "Bibliografia:" / (recovery norm text stays) ---
$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = "Média aritmética da nota final obtida pelo aluno durante o semestre e da nota obtida na Prova de Recuperação."
$ws.Range("C21").Value = "Média aritmética da nota final obtida pelo aluno durante o semestre e da nota obtida na Prova de Recuperação."

# --- Row 22: becomes just "Requisitos:" (B/C cleared entirely) ---
$ws.Range("A22").Value = "Requisitos:"
$ws.Range("B22:C22").Clear()

# --- Row 23: A cleared, B/C become the requirement detail text (newly populated) ---
$ws.Range("A23").Clear()
$ws.Range("B19").Copy()
$ws.Range("B23").PasteSpecial($xlPasteFormats)
$ws.Range("C19").Copy()
$ws.Range("C23").PasteSpecial($xlPasteFormats)
$ws.Range("B23").Value = "LOQ4205 -  Sistemas Produtivos II  (Requisito fraco)" + [char]10
$ws.Range("C23").Value = "LOQ4205 -  Sistemas Produtivos II  (Requisito fraco)" + [char]10

# --- Row 24: no longer exists in the target sheet; remove it entirely ---
$ws.Rows.Item(24).Delete()

# --- Fix up row heights to match the final target layout ---
$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(14).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 120
$ws.Rows.Item(16).RowHeight = 120
# rows 17 and 22 go back to the sheet's default (un-customized) row height
$ws.Rows.Item(17).AutoFit()
$ws.Rows.Item(18).RowHeight = 60
$ws.Rows.Item(19).RowHeight = 60
$ws.Rows.Item(20).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 120
$ws.Rows.Item(22).AutoFit()
$ws.Rows.Item(23).RowHeight = 30
